$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: Geometry Instancing - clear milestone marker (Cleaned Up Some Code)
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""

# Row 23: Infinite Sky Box - mark milestone II complete (Added Sky Box)
$ws.Range("E23").Value = "II"
$ws.Range("F23").Value = "X"

# Row 32: Functional per pixel spot light - mark milestone II complete (Added Spot Light)
$ws.Range("E32").Value = "II"
$ws.Range("F32").Value = "X"

# Row 36: Dynamic change in position & direction of spot lighting - mark milestone II complete
$ws.Range("E36").Value = "II"
$ws.Range("F36").Value = "X"

# Rows 90-91: mark Milestone II complete for Effective Use of GIT / Confidence Confirmed (Changed plane)
$ws.Range("D90").Value = "X"
$ws.Range("D91").Value = "X"

# Update the view: scroll so column C is left-most visible and select E67
$ws.Range("E67").Select()
$excel.ActiveWindow.ScrollColumn = 3
